$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 20,10
$arr[0,0] = -18.23226632836297
$arr[0,1] = -18.23226632836297
$arr[0,2] = -18.23226632836297
$arr[0,3] = -18.23226632836297
$arr[0,4] = -18.23226632836297
$arr[0,5] = -18.23226632836297
$arr[0,6] = -18.23226632836297
$arr[0,7] = -18.23226632836297
$arr[0,8] = -18.23226632836297
$arr[0,9] = -18.23226632836297
$arr[1,0] = -18.23226632836297
$arr[1,1] = -18.23226632836297
$arr[1,2] = -18.23226632836297
$arr[1,3] = -18.23226632836297
$arr[1,4] = -18.23226632836297
$arr[1,5] = -18.23226632836297
$arr[1,6] = -18.23226632836297
$arr[1,7] = 1.257199478593597
$arr[1,8] = -18.23226632836297
$arr[1,9] = -18.23226632836297
$arr[2,0] = -18.23226632836297
$arr[2,1] = -18.23226632836297
$arr[2,2] = 1.896590187191784
$arr[2,3] = -18.23226632836297
$arr[2,4] = 3.54663415231373
$arr[2,5] = -18.23226632836297
$arr[2,6] = 1.468974130412395
$arr[2,7] = -18.23226632836297
$arr[2,8] = 0.8031583697578338
$arr[2,9] = -18.23226632836297
$arr[3,0] = -18.23226632836297
$arr[3,1] = -18.23226632836297
$arr[3,2] = -18.23226632836297
$arr[3,3] = -18.23226632836297
$arr[3,4] = -18.23226632836297
$arr[3,5] = 2.860777850133727
$arr[3,6] = -18.23226632836297
$arr[3,7] = -18.23226632836297
$arr[3,8] = -18.23226632836297
$arr[3,9] = -18.23226632836297
$arr[4,0] = -18.23226632836297
$arr[4,1] = -18.23226632836297
$arr[4,2] = -18.23226632836297
$arr[4,3] = -18.23226632836297
$arr[4,4] = -18.23226632836297
$arr[4,5] = -18.23226632836297
$arr[4,6] = -18.23226632836297
$arr[4,7] = -18.23226632836297
$arr[4,8] = -18.23226632836297
$arr[4,9] = -18.23226632836297
$arr[5,0] = 2.370650874172485
$arr[5,1] = -18.23226632836297
$arr[5,2] = -18.23226632836297
$arr[5,3] = -18.23226632836297
$arr[5,4] = -18.23226632836297
$arr[5,5] = -18.23226632836297
$arr[5,6] = -18.23226632836297
$arr[5,7] = -18.23226632836297
$arr[5,8] = -18.23226632836297
$arr[5,9] = -18.23226632836297
$arr[6,0] = -18.23226632836297
$arr[6,1] = -18.23226632836297
$arr[6,2] = -18.23226632836297
$arr[6,3] = 1.901740023772815
$arr[6,4] = -18.23226632836297
$arr[6,5] = -18.23226632836297
$arr[6,6] = -18.23226632836297
$arr[6,7] = -18.23226632836297
$arr[6,8] = -18.23226632836297
$arr[6,9] = -18.23226632836297
$arr[7,0] = 3.890271331751803
$arr[7,1] = -18.23226632836297
$arr[7,2] = -18.23226632836297
$arr[7,3] = -18.23226632836297
$arr[7,4] = -18.23226632836297
$arr[7,5] = -18.23226632836297
$arr[7,6] = -18.23226632836297
$arr[7,7] = -18.23226632836297
$arr[7,8] = -18.23226632836297
$arr[7,9] = -18.23226632836297
$arr[8,0] = -18.23226632836297
$arr[8,1] = -18.23226632836297
$arr[8,2] = -18.23226632836297
$arr[8,3] = -18.23226632836297
$arr[8,4] = -18.23226632836297
$arr[8,5] = -18.23226632836297
$arr[8,6] = -18.23226632836297
$arr[8,7] = 1.711543365844584
$arr[8,8] = -18.23226632836297
$arr[8,9] = 2.21180322613731
$arr[9,0] = -18.23226632836297
$arr[9,1] = -18.23226632836297
$arr[9,2] = -18.23226632836297
$arr[9,3] = 3.002165801035235
$arr[9,4] = -18.23226632836297
$arr[9,5] = 2.845332420499925
$arr[9,6] = -18.23226632836297
$arr[9,7] = -18.23226632836297
$arr[9,8] = -18.23226632836297
$arr[9,9] = 1.947376749947417
$arr[10,0] = -18.23226632836297
$arr[10,1] = -18.23226632836297
$arr[10,2] = -18.23226632836297
$arr[10,3] = -18.23226632836297
$arr[10,4] = -18.23226632836297
$arr[10,5] = -18.23226632836297
$arr[10,6] = -18.23226632836297
$arr[10,7] = -18.23226632836297
$arr[10,8] = -18.23226632836297
$arr[10,9] = -18.23226632836297
$arr[11,0] = -18.23226632836297
$arr[11,1] = -18.23226632836297
$arr[11,2] = -18.23226632836297
$arr[11,3] = 2.496033710114683
$arr[11,4] = -18.23226632836297
$arr[11,5] = -18.23226632836297
$arr[11,6] = -18.23226632836297
$arr[11,7] = -18.23226632836297
$arr[11,8] = 1.716798045726078
$arr[11,9] = 1.825589257735758
$arr[12,0] = -18.23226632836297
$arr[12,1] = -18.23226632836297
$arr[12,2] = 1.342513806013164
$arr[12,3] = -18.23226632836297
$arr[12,4] = -18.23226632836297
$arr[12,5] = -18.23226632836297
$arr[12,6] = -18.23226632836297
$arr[12,7] = -18.23226632836297
$arr[12,8] = -18.23226632836297
$arr[12,9] = 1.947060873269506
$arr[13,0] = -18.23226632836297
$arr[13,1] = -18.23226632836297
$arr[13,2] = 1.63679880492824
$arr[13,3] = -18.23226632836297
$arr[13,4] = -18.23226632836297
$arr[13,5] = -18.23226632836297
$arr[13,6] = -18.23226632836297
$arr[13,7] = -18.23226632836297
$arr[13,8] = -18.23226632836297
$arr[13,9] = -18.23226632836297
$arr[14,0] = -18.23226632836297
$arr[14,1] = -18.23226632836297
$arr[14,2] = -18.23226632836297
$arr[14,3] = -18.23226632836297
$arr[14,4] = -18.23226632836297
$arr[14,5] = -18.23226632836297
$arr[14,6] = -18.23226632836297
$arr[14,7] = -18.23226632836297
$arr[14,8] = 1.932774362254625
$arr[14,9] = -18.23226632836297
$arr[15,0] = -18.23226632836297
$arr[15,1] = -18.23226632836297
$arr[15,2] = 2.055627503760666
$arr[15,3] = -18.23226632836297
$arr[15,4] = -18.23226632836297
$arr[15,5] = -18.23226632836297
$arr[15,6] = 2.03919327992115
$arr[15,7] = 2.091664730081163
$arr[15,8] = 2.543383333041013
$arr[15,9] = -18.23226632836297
$arr[16,0] = -18.23226632836297
$arr[16,1] = -18.23226632836297
$arr[16,2] = -18.23226632836297
$arr[16,3] = -18.23226632836297
$arr[16,4] = -18.23226632836297
$arr[16,5] = -18.23226632836297
$arr[16,6] = 1.986917557339842
$arr[16,7] = 2.073480387055346
$arr[16,8] = 2.411573130842471
$arr[16,9] = -18.23226632836297
$arr[17,0] = -18.23226632836297
$arr[17,1] = -18.23226632836297
$arr[17,2] = 1.969102625294671
$arr[17,3] = -18.23226632836297
$arr[17,4] = -18.23226632836297
$arr[17,5] = -18.23226632836297
$arr[17,6] = 1.633873090270733
$arr[17,7] = 1.838170433201685
$arr[17,8] = -18.23226632836297
$arr[17,9] = -18.23226632836297
$arr[18,0] = -18.23226632836297
$arr[18,1] = 4.321923644070135
$arr[18,2] = 1.355278644210777
$arr[18,3] = -18.23226632836297
$arr[18,4] = 3.055638209876565
$arr[18,5] = -18.23226632836297
$arr[18,6] = 1.696791210310593
$arr[18,7] = 1.193812065101366
$arr[18,8] = -18.23226632836297
$arr[18,9] = 2.03923887419146
$arr[19,0] = -18.23226632836297
$arr[19,1] = -18.23226632836297
$arr[19,2] = -18.23226632836297
$arr[19,3] = 1.384033023040793
$arr[19,4] = -18.23226632836297
$arr[19,5] = 2.472269751808247
$arr[19,6] = 1.491809472522965
$arr[19,7] = -18.23226632836297
$arr[19,8] = -18.23226632836297
$arr[19,9] = -18.23226632836297
$ws.Range("B2:K21").Value = $arr
Write-Output "done"
